$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H, cloning the formatting of the existing header cells (B1:G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New "Save" flag column values (row 2 through row 8)
$saveValues = @(0, 0, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
